# "understanfd the varience of data"
# Adds a second worked example (3 students, "Classe A" / "Classe B") below the
# existing note-1/note-2 example: a small step-by-step population-variance
# computation, shown three times side by side (raw numbers, "xi - avg",
# "(xi - avg)^2"), with running totals, the mean of the squared deviations,
# and a final cross-check against VAR.P().

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yellow = 65535      # RGB(255,255,0)
$blueTint = 15983578 # RGB(0xDA,0xE3,0xF3)  ~ Accent1 (theme4), Lighter 80%
$orangeTint = 11389944 # RGB(0xF8,0xCB,0xAD) ~ Accent2 (theme5), Lighter 60%

function Style-Header($rng, $fill, $halign, $border) {
    $rng.Font.Bold = $true
    if ($fill -ne $null) { $rng.Interior.Color = $fill }
    $rng.WrapText = $true
    $rng.VerticalAlignment = -4108
    $rng.HorizontalAlignment = $halign
    if ($border) { $rng.Borders.LineStyle = 1 }
}

function Style-Data($rng, $fill, $halign, $border) {
    if ($fill -ne $null) { $rng.Interior.Color = $fill }
    $rng.WrapText = $true
    $rng.VerticalAlignment = -4108
    if ($halign -ne $null) { $rng.HorizontalAlignment = $halign }
    if ($border) { $rng.Borders.LineStyle = 1 }
}

function Style-Total($rng) {
    $rng.Font.Bold = $true
    $rng.Interior.Color = $yellow
}

# ---------------------------------------------------------------------
# Row 40 banner cells (merged, yellow, bold, boxed, centered)
# ---------------------------------------------------------------------
$banner1 = $ws.Range("F40:G40")
$banner1.Merge()
$banner1.Value = "xi - avg"
$banner1.Font.Bold = $true
$banner1.Interior.Color = $yellow
$banner1.HorizontalAlignment = -4108
$banner1.Borders.LineStyle = 1

$banner2 = $ws.Range("J40:K40")
$banner2.Merge()
$banner2.Value = "(xi - avg)2"
$banner2.Font.Bold = $true
$banner2.Interior.Color = $yellow
$banner2.HorizontalAlignment = -4108
$banner2.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# Row 41 headers - three parallel tables (A:C, E:G, I:K)
# ---------------------------------------------------------------------

# --- Table 1 (A:C) : no outer borders ---
Style-Header $ws.Range("A41") $null -4108 $false
$ws.Range("A41").Value = "Élève"

Style-Header $ws.Range("B41") $blueTint -4108 $false
$ws.Range("B41").Value = "Classe A"

Style-Header $ws.Range("C41") $orangeTint -4108 $false
$ws.Range("C41").Value = "Classe B"

# --- Table 2 (E:G) : boxed ---
Style-Header $ws.Range("E41") $null -4131 $true
$ws.Range("E41").Value = "Élève"

Style-Header $ws.Range("F41") $blueTint -4108 $true
$ws.Range("F41").Value = "Classe A"

Style-Header $ws.Range("G41") $orangeTint -4108 $true
$ws.Range("G41").Value = "Classe B"

# --- Table 3 (I:K) : boxed ---
Style-Header $ws.Range("I41") $null -4131 $true
$ws.Range("I41").Value = "Élève"

Style-Header $ws.Range("J41") $blueTint -4108 $true
$ws.Range("J41").Value = "Classe A"

Style-Header $ws.Range("K41") $orangeTint -4108 $true
$ws.Range("K41").Value = "Classe B"

# ---------------------------------------------------------------------
# Rows 42:44 - student index / raw data / deviations / squared deviations
# ---------------------------------------------------------------------

# Student index column (no border in table1, bordered in table2/3)
Style-Data $ws.Range("A42:A44") $null $null $false
$ws.Range("A42").Value = 1
$ws.Range("A43").Value = 2
$ws.Range("A44").Value = 3

Style-Data $ws.Range("E42:E44") $null -4131 $true
$ws.Range("E42").Value = 1
$ws.Range("E43").Value = 2
$ws.Range("E44").Value = 3

Style-Data $ws.Range("I42:I44") $null -4131 $true
$ws.Range("I42").Value = 1
$ws.Range("I43").Value = 2
$ws.Range("I44").Value = 3

# Raw data, Classe A / Classe B (table 1)
Style-Data $ws.Range("B42:B44") $blueTint $null $false
$ws.Range("B42").Value = 10
$ws.Range("B43").Value = 10
$ws.Range("B44").Value = 10

Style-Data $ws.Range("C42:C44") $orangeTint $null $false
$ws.Range("C42").Value = 6
$ws.Range("C43").Value = 10
$ws.Range("C44").Value = 14

# xi - avg (table 2), first row entered individually then filled down
# (matches how Excel records a lone formula cell plus a shared range)
Style-Data $ws.Range("F42:F44") $blueTint $null $true
$ws.Range("F42").Formula = "=B42-`$B`$46"
$ws.Range("F43:F44").Formula = "=B43-`$B`$46"

Style-Data $ws.Range("G42:G44") $orangeTint $null $true
$ws.Range("G42").Formula = "=C42-`$C`$46"
$ws.Range("G43:G44").Formula = "=C43-`$C`$46"

# (xi - avg)^2 (table 3)
Style-Data $ws.Range("J42:J44") $blueTint $null $true
$ws.Range("J42").Formula = "=POWER(F42,2)"
$ws.Range("J43:J44").Formula = "=POWER(F43,2)"

Style-Data $ws.Range("K42:K44") $orangeTint $null $true
$ws.Range("K42").Formula = "=POWER(G42,2)"
$ws.Range("K43:K44").Formula = "=POWER(G43,2)"

# ---------------------------------------------------------------------
# Separator rows 45 / 50 (thick-bottom spacer rows above the totals)
# ---------------------------------------------------------------------
$ws.Rows.Item(45).RowHeight = 15
$ws.Rows.Item(46).RowHeight = 15
$ws.Rows.Item(50).RowHeight = 15
$ws.Rows.Item(51).RowHeight = 15

# ---------------------------------------------------------------------
# Row 46 - Moyenne / Somme / Somme totals (medium boxed, yellow)
# ---------------------------------------------------------------------
$ws.Range("A46").Value = "Moyenne"
$ws.Range("B46").Formula = "=AVERAGE(B42:B44)"
$ws.Range("C46").Formula = "=AVERAGE(C42:C44)"
Style-Total $ws.Range("A46:C46")
$ws.Range("A46:C46").BorderAround($null, -4138)

$ws.Range("E46").Value = "Somme"
$ws.Range("F46").Formula = "=SUM(F42:F44)"
$ws.Range("G46").Formula = "=SUM(G42:G44)"
Style-Total $ws.Range("E46:G46")
$ws.Range("E46:G46").BorderAround($null, -4138)

$ws.Range("I46").Value = "Somme"
$ws.Range("J46").Formula = "=SUM(J42:J44)"
$ws.Range("K46").Formula = "=SUM(K42:K44)"
Style-Total $ws.Range("I46:K46")
$ws.Range("I46:K46").BorderAround($null, -4138)

# ---------------------------------------------------------------------
# Row 51 - moy (mean of the squared deviations => population variance)
# ---------------------------------------------------------------------
$ws.Range("I51").Value = "moy"
$ws.Range("J51").Formula = "=J46/3"
$ws.Range("K51").Formula = "=K46/3"
Style-Total $ws.Range("I51:K51")
$ws.Range("I51:K51").BorderAround($null, -4138)

# ---------------------------------------------------------------------
# Row 53 - cross-check against the built-in VAR.P function
# ---------------------------------------------------------------------
$ws.Range("J53").Value = "varPop"
$ws.Range("K53").Formula = "=VAR.P(C42:C44)"
$ws.Range("J53:L53").Interior.Color = $orangeTint

# ---------------------------------------------------------------------
# View bookkeeping - mirror where the author left the selection/scroll
# ---------------------------------------------------------------------
$ws.Range("S37").Select()
